$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3400
$ws1.Range("F5").Value = 1592
$ws1.Range("F7").Value = 328

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3400
$ws4.Range("F5").Value = 1592
$ws4.Range("F8").Value = 328
